# Applies the semantic changes described in the commit:
#   "Fix ExcelIO: bool & float write. Fix exceljs: table position write,
#    cell value type read and write, confusing naming. ..."
#
# Concretely (per the xlsx canonical-XML diff):
#   - On every sheet, the "Numbers" column's last data cell (row 5 on
#     WithTable/Tableless, row 8 on WithTable_Duplicate) changes from the
#     integer 4 to the float 4.269, formatted with number format "0.00".
#   - The previously selected cell on every sheet moves to that same cell
#     (A5 / A5 / B8).
#   - The active sheet changes from "WithTable_Duplicate" (index 2) to
#     "Tableless" (index 1), which also moves which sheet carries
#     tabSelected="1".

$wb = $excel.ActiveWorkbook

# --- WithTable sheet: A5 becomes 4.269 with a 0.00 number format -----------
$wsWithTable = $wb.Worksheets.Item("WithTable")
$wsWithTable.Range("A5").Value = 4.269
$wsWithTable.Range("A5").NumberFormat = "0.00"
$wsWithTable.Range("A5").Select()

# --- WithTable_Duplicate sheet: B8 becomes 4.269 with a 0.00 number format -
$wsDuplicate = $wb.Worksheets.Item("WithTable_Duplicate")
$wsDuplicate.Range("B8").Value = 4.269
$wsDuplicate.Range("B8").NumberFormat = "0.00"
$wsDuplicate.Range("B8").Select()

# --- Tableless sheet: A5 becomes 4.269 with a 0.00 number format -----------
# Selected/activated last so it ends up as the workbook's active tab,
# matching the target's activeTab/tabSelected state.
$wsTableless = $wb.Worksheets.Item("Tableless")
$wsTableless.Range("A5").Value = 4.269
$wsTableless.Range("A5").NumberFormat = "0.00"
$wsTableless.Range("A5").Select()
